# Fixed naive component forecaster bug - Presentation state 11.02.
# Update forecast error metrics (ME, MAE, MSE, RMSE, SE, N) for quarters Q1..Q9 (rows 2-11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.3291291692495283, 1.812906911408057, 27.04923173616717, 5.200887591187409, 5.247817881751421, 46)
    3  = @(-0.9495326209681835, 2.023872778234117, 43.75773521907136, 6.614962979418053, 6.620432449373355, 45)
    4  = @(-1.78880049929551,   3.145904294961889, 171.4597706064012, 13.094264798239,   13.12147041243131, 44)
    5  = @(-4.213288228997171,  5.762860054837126, 837.7929586707833, 28.94465336933202, 28.97526445982997, 43)
    6  = @(-8.901962788184429,  10.44294745540641, 3576.380529151133, 59.80284716592625, 59.85341742634296, 42)
    7  = @(-18.90835081736334,  20.50494998719962, 14831.91797591911, 121.786362027606,  121.8041544916215, 41)
    8  = @(-39.82066791606469,  41.45267336896583, 64161.89026247886, 253.3019744543632, 253.3391456874329, 40)
    9  = @(-84.10188452280053,  85.72607652365316, 276745.0078080686, 526.0655926859963, 526.087929219649,  39)
    10 = @(-176.9409422996708,  178.6784928052071, 1193366.26590965,  1092.413047299258, 1092.458304185648, 38)
    11 = @(-372.8822845639866,  374.5397916143799, 5147954.627584639, 2268.910449441458, 2268.931447530792, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
